$d = $word.ActiveDocument

# Locate the end of the paragraph that currently ends the "搭建微服务架构" section:
# "Docker 轻量级的特点使得它很适合用于部署、维护、组合微服务。"
$rng = $d.Content
$found = $rng.Find.Execute("轻量级的特点使得它很适合用于部署、维护、组合微服务。", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph text"
}

# Collapse to the end of the found text (end of that paragraph's content).
$rng.Collapse(0)

# Insert an empty paragraph mark, then the new "五、镜像与容器" section content,
# each new sentence/paragraph separated by a paragraph mark (`r`).
$newText = "`r`r五、镜像与容器`r镜像是一种静态的结构，可以看成面向对象里面的类，而容器是镜像的一个实例。`r镜像包含着容器运行时所需要的代码以及其它组件，它是一种分层结构，每一层都是只读的（"
$rng.InsertAfter($newText)

$rng.Collapse(0)
$rng.InsertAfter("read-only layers")

$rng.Collapse(0)
$rng.InsertAfter("）。构建镜像时，会一层一层构建，前一层是后一层的基础。镜像的这种分层存储结构很适合镜像的复用以及定制。`r构建容器时，通过在镜像的基础上添加一个可写层（")

$rng.Collapse(0)
$rng.InsertAfter("writable layer")

$rng.Collapse(0)
$rng.InsertAfter("），用来保存着容器运行过程中的修改。")
